# This script applies updated betting-odds values (and one kickoff-time correction)
# to Sheet1 of the FlashScore weekly games workbook, matching the target diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 2
$ws.Range("O2").Value = 1.62   # was 1.57
$ws.Range("P2").Value = 2.2   # was 2.25
$ws.Range("Q2").Value = 2.21   # was 2.18
$ws.Range("R2").Value = 1.67   # was 1.69
$ws.Range("U2").Value = 5.2   # was 5

# Row 3
$ws.Range("I3").Value = 5.75   # was 6
$ws.Range("K3").Value = 2   # was 1.95
$ws.Range("Q3").Value = 1.93   # was 1.95
$ws.Range("R3").Value = 1.93   # was 1.9
$ws.Range("AA3").Value = 2.38   # was 2.5
$ws.Range("AB3").Value = 1.53   # was 1.5
$ws.Range("AC3").Value = 5   # was 4.75
$ws.Range("AO3").Value = 26   # was 29
$ws.Range("AP3").Value = 19   # was 21

# Row 4
$ws.Range("AC4").Value = 5.5   # was 5

# Row 5
$ws.Range("S5").Value = 2.15   # was 2.2
$ws.Range("T5").Value = 1.67   # was 1.65
$ws.Range("U5").Value = 3.35   # was 3.5
$ws.Range("V5").Value = 1.33   # was 1.31

# Row 6
$ws.Range("G6").Value = 3.8   # was 3.9
$ws.Range("I6").Value = 2.2   # was 2.15
$ws.Range("M6").Value = 1.1   # was 1.08
$ws.Range("N6").Value = 7   # was 8
$ws.Range("O6").Value = 1.44   # was 1.4
$ws.Range("P6").Value = 2.63   # was 2.75
$ws.Range("Q6").Value = 1.83   # was 1.74
$ws.Range("R6").Value = 2.03   # was 2.1
$ws.Range("S6").Value = 2.4   # was 2.3
$ws.Range("T6").Value = 1.53   # was 1.6
$ws.Range("U6").Value = 3.8   # was 3.6
$ws.Range("V6").Value = 1.27   # was 1.3
$ws.Range("W6").Value = 4.5   # was 4.33
$ws.Range("X6").Value = 1.18   # was 1.2
$ws.Range("Y6").Value = 1.53   # was 1.5
$ws.Range("Z6").Value = 2.38   # was 2.5
$ws.Range("AC6").Value = 9   # was 9.5
$ws.Range("AJ6").Value = 5.5   # was 6
$ws.Range("AO6").Value = 9.5   # was 9
$ws.Range("AR6").Value = 21   # was 19

# Row 7
$ws.Range("H7").Value = 5.8   # was 5.9
$ws.Range("L7").Value = 7.2   # was 7.3
$ws.Range("AB7").Value = 1.83   # was 1.82
$ws.Range("AC7").Value = 10.25   # was 10

# Row 8
$ws.Range("C8").Value = "21:30"   # was "20:30"
$ws.Range("N8").Value = 8   # was 7.5
$ws.Range("S8").Value = 2.25   # was 2.3
$ws.Range("T8").Value = 1.62   # was 1.6
$ws.Range("W8").Value = 4   # was 4.33
$ws.Range("X8").Value = 1.22   # was 1.2

# Row 10
$ws.Range("G10").Value = 2.9   # was 2.87
$ws.Range("H10").Value = 2.72   # was 2.77
$ws.Range("I10").Value = 2.72   # was 2.7
$ws.Range("J10").Value = 3.6   # was 3.55
$ws.Range("K10").Value = 1.83   # was 1.88
$ws.Range("L10").Value = 3.45   # was 3.35
$ws.Range("N10").Value = 5   # was 5.2
$ws.Range("O10").Value = 1.55   # was 1.53
$ws.Range("P10").Value = 2.3   # was 2.32
$ws.Range("S10").Value = 2.62   # was 2.55
$ws.Range("T10").Value = 1.44   # was 1.45
$ws.Range("W10").Value = 4.75   # was 4.6
$ws.Range("X10").Value = 1.15   # was 1.16
$ws.Range("Y10").Value = 1.6   # was 1.55
$ws.Range("Z10").Value = 2.22   # was 2.3
$ws.Range("AA10").Value = 2.1   # was 2.05
$ws.Range("AB10").Value = 1.65   # was 1.7
$ws.Range("AC10").Value = 6.6   # was 6.7
$ws.Range("AG10").Value = 32   # was 30
$ws.Range("AI10").Value = 5   # was 5.2
$ws.Range("AJ10").Value = 5.4   # was 5.5
$ws.Range("AN10").Value = 6.3   # was 6.4
$ws.Range("AP10").Value = 10.75   # was 10.5
$ws.Range("AR10").Value = 30   # was 28
$ws.Range("AS10").Value = 50   # was 45

# Row 11
$ws.Range("G11").Value = 4.85   # was 4.8
$ws.Range("H11").Value = 3.15   # was 3.2
$ws.Range("J11").Value = 5.2   # was 5
$ws.Range("L11").Value = 2.35   # was 2.37
$ws.Range("Z11").Value = 2.55   # was 2.57
$ws.Range("AA11").Value = 1.98   # was 2
$ws.Range("AB11").Value = 1.75   # was 1.72
$ws.Range("AC11").Value = 10.5   # was 11.25
$ws.Range("AD11").Value = 26   # was 27
$ws.Range("AE11").Value = 16   # was 15.5
$ws.Range("AF11").Value = 100   # was 90
$ws.Range("AG11").Value = 60   # was 55
$ws.Range("AH11").Value = 65   # was 60
$ws.Range("AL11").Value = 100   # was 90
$ws.Range("AN11").Value = 5.8   # was 5.7
$ws.Range("AO11").Value = 7.8   # was 7.5
$ws.Range("AR11").Value = 15.5   # was 16

# Row 12
$ws.Range("I12").Value = 12.5   # was 11.75
$ws.Range("K12").Value = 2.32   # was 2.35
$ws.Range("L12").Value = 10   # was 9.5
$ws.Range("N12").Value = 7   # was 7.1
$ws.Range("O12").Value = 1.33   # was 1.32
$ws.Range("P12").Value = 3.05   # was 3.1
$ws.Range("S12").Value = 1.98   # was 1.95
$ws.Range("W12").Value = 3.3   # was 3.25
$ws.Range("Y12").Value = 1.39   # was 1.37
$ws.Range("Z12").Value = 2.77   # was 2.85
$ws.Range("AA12").Value = 2.57   # was 2.55
$ws.Range("AB12").Value = 1.44   # was 1.45
$ws.Range("AG12").Value = 13.5   # was 13
$ws.Range("AI12").Value = 7   # was 7.1
$ws.Range("AL12").Value = 250   # was 200
$ws.Range("AP12").Value = 40   # was 37
$ws.Range("AQ12").Value = 500   # was 450
$ws.Range("AR12").Value = 250   # was 200

# Row 14
$ws.Range("I14").Value = 3.7   # was 3.75
$ws.Range("J14").Value = 2.88   # was 2.75
$ws.Range("L14").Value = 4.33   # was 4.5
$ws.Range("N14").Value = 7.5   # was 8
$ws.Range("AF14").Value = 19   # was 17
$ws.Range("AI14").Value = 7.5   # was 8
$ws.Range("AL14").Value = 67   # was 51
$ws.Range("AM14").Value = 451   # was 401
$ws.Range("AN14").Value = 9   # was 9.5
$ws.Range("AO14").Value = 17   # was 19

# Row 15
$ws.Range("I15").Value = 3.1   # was 3.2
$ws.Range("J15").Value = 3   # was 2.88
$ws.Range("K15").Value = 2.1   # was 2.05
$ws.Range("L15").Value = 3.75   # was 4
$ws.Range("M15").Value = 1.06   # was 1.07
$ws.Range("N15").Value = 10   # was 9
$ws.Range("O15").Value = 1.33   # was 1.36
$ws.Range("P15").Value = 3.4   # was 3.2
$ws.Range("S15").Value = 2.05   # was 2.1
$ws.Range("T15").Value = 1.8   # was 1.73
$ws.Range("AA15").Value = 1.8   # was 1.83
$ws.Range("AB15").Value = 1.91   # was 1.83
$ws.Range("AC15").Value = 7.5   # was 7
$ws.Range("AD15").Value = 11   # was 10
$ws.Range("AE15").Value = 9   # was 9.5
$ws.Range("AI15").Value = 10   # was 9
$ws.Range("AM15").Value = 251   # was 301
$ws.Range("AN15").Value = 9.5   # was 9
$ws.Range("AP15").Value = 11   # was 12

# Row 17
$ws.Range("M17").Value = 1.06   # was 1.07
$ws.Range("N17").Value = 10   # was 9

# Row 18
$ws.Range("S18").Value = 1.73   # was 1.7
$ws.Range("T18").Value = 2.08   # was 2.1

# Row 19
$ws.Range("M19").Value = 1.06   # was 1.07
$ws.Range("N19").Value = 10   # was 9
$ws.Range("W19").Value = 3.5   # was 3.75
$ws.Range("X19").Value = 1.29   # was 1.25

# Row 20
$ws.Range("M20").Value = 1.05   # was 1.06
$ws.Range("N20").Value = 11   # was 9.5
$ws.Range("O20").Value = 1.29   # was 1.3
$ws.Range("P20").Value = 3.5   # was 3.4
$ws.Range("S20").Value = 2   # was 2.03
$ws.Range("T20").Value = 1.85   # was 1.83
$ws.Range("W20").Value = 3.4   # was 3.5
$ws.Range("X20").Value = 1.3   # was 1.29

# Row 21
$ws.Range("M21").Value = 1.03   # was 1.04
$ws.Range("N21").Value = 15   # was 13

# Row 22
$ws.Range("S22").Value = 2.2   # was 2.15
$ws.Range("T22").Value = 1.65   # was 1.67

# Row 23
$ws.Range("G23").Value = 1.91   # was 1.9
$ws.Range("I23").Value = 3.7   # was 3.75
$ws.Range("AC23").Value = 7   # was 6.5
$ws.Range("AF23").Value = 17   # was 15

# Row 24
$ws.Range("S24").Value = 1.85   # was 1.88
$ws.Range("T24").Value = 2   # was 1.98

# Row 26
$ws.Range("K26").Value = 2.1   # was 2.2
$ws.Range("M26").Value = 1.06   # was 1.05
$ws.Range("N26").Value = 10   # was 11
$ws.Range("AC26").Value = 10   # was 11
$ws.Range("AI26").Value = 10   # was 11
$ws.Range("AN26").Value = 8   # was 8.5

